$d = $word.ActiveDocument
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd("`r")
}

# ---------------------------------------------------------------------
# 1) Expand the "Wireless Sensor Networks" section with new content and
#    add a new "What is the target application of your project?"
#    sub-section, right before the existing "Concept description" /
#    "Block diagram..." / "What is the main application for your
#    prototype?" paragraphs (those are kept, just pushed down).
# ---------------------------------------------------------------------

# Locate the (empty) paragraph that only holds the stray "_GoBack"
# bookmark -- it immediately follows the "Wireless Sensor Networks"
# heading. Remove it; the bookmark itself will be re-created further
# down, in its new resting place.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "") {
        $prev = $p.Previous()
        if ($prev -ne $null -and (Get-ParaText $prev) -eq "Wireless Sensor Networks ") {
            $targetPara = $p
            break
        }
    }
}
$targetPara.Range.Delete()

# Find the run of paragraphs that now starts right where the deleted
# paragraph used to be: "Concept description", "Block diagram of your
# target application.", "What is the main application for your
# prototype?".
$pConcept = $null
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "Concept description") {
        $pConcept = $p
        break
    }
}
$pBlockDiagram = $pConcept.Next()
$pMainApplication = $pBlockDiagram.Next()

$rng = $d.Range($pConcept.Range.Start, $pMainApplication.Range.End)

$newXml = @"
<w:p $w>
  <w:r><w:t>Wireless sensor networks are collections of small sensors that can sense and transmit data wirelessly, such as temperature or object detection. The color sensor and ultrasonic sensor are a part of this network in our project. They transmit data wirelessly (via an Arduino and a Raspberry Pi) to the main controllers. This facilitates wireless monitoring and control of the robotic arm and conveyor belt.</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Heading2"/></w:pPr>
  <w:r><w:t>What is the target application of your project?</w:t></w:r>
</w:p>
<w:p $w>
  <w:r><w:t xml:space="preserve">Our research aims to implement an automated system for sorting small colored cubes in an industrial environment. The cubes are moved by the conveyor belt; when a cube is identified, the robotic arm detects it and sorts the cubes according to color. This system can help industries sort items quickly and efficiently without needing human </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>intervention.</w:t></w:r>
</w:p>
<w:p $w>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p $w>
  <w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
  <w:r><w:t>Concept</w:t></w:r>
  <w:r><w:t xml:space="preserve"> description</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr><w:jc w:val="both"/><w:rPr><w:i/></w:rPr></w:pPr>
  <w:r><w:rPr><w:i/></w:rPr><w:t>Block diagram of your target application.</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr><w:jc w:val="both"/><w:rPr><w:i/></w:rPr></w:pPr>
  <w:r><w:rPr><w:i/></w:rPr><w:t>What is the main application for</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> your prototype?</w:t></w:r>
</w:p>
"@

$rng.InsertXML($newXml)

# ---------------------------------------------------------------------
# 2) Move the "lastRenderedPageBreak" marker: it used to sit on the
#    "..." bullet (Technologies section); it now belongs on the
#    "Breakdown: " run (Project/Team management section).
# ---------------------------------------------------------------------

$pBreakdown = $null
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "Breakdown: How you managed your tasks?") {
        $pBreakdown = $p
        break
    }
}

$breakdownXml = @"
<w:p $w>
  <w:pPr><w:jc w:val="both"/><w:rPr><w:i/></w:rPr></w:pPr>
  <w:r><w:rPr><w:i/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Breakdown: </w:t></w:r>
  <w:r><w:rPr><w:i/></w:rPr><w:t>How you managed your tasks?</w:t></w:r>
</w:p>
"@

$pBreakdown.Range.InsertXML($breakdownXml)

$pDots = $null
foreach ($p in $d.Paragraphs) {
    if ((Get-ParaText $p) -eq "...") {
        $pDots = $p
        break
    }
}

$dotsXml = @"
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="29"/></w:numPr>
    <w:ind w:left="773"/>
    <w:jc w:val="both"/>
    <w:rPr><w:i/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:i/></w:rPr><w:t>...</w:t></w:r>
</w:p>
"@

$pDots.Range.InsertXML($dotsXml)
